$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Set header values in the newly inserted row 1
$ws.Range("A1").Value = "category"
$ws.Range("B1").Value = "treated"
$ws.Range("C1").Value = "control"
$ws.Range("D1").Value = "pvalue"
